# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (right before "总计") holding the
#   per-fund holdings detail for the quarter.
# - Insert a new leading row into the "总计" (totals) sheet summarising
#   the 2022-Q1 holdings, shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Scratch sheet used only to stage values that must be written as TEXT
# (e.g. fund codes with leading zeros, or numeric-looking ratios that
# are stored as text in this workbook) without leaving a NumberFormat/
# style footprint on the real cells. Range.PasteSpecial(xlPasteValues)
# copies the cell's type+value but not its formatting, so the
# destination cell ends up exactly like a plain inlineStr text cell.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Add()
$scratch.Name = "__scratch__"
$scratch.Range("A1:G1").NumberFormat = "@"
$scratch.Range("A1").Value = "005416"
$scratch.Range("B1").Value = "7.95"
$scratch.Range("C1").Value = "37.81"
$scratch.Range("D1").Value = "1.24"
$scratch.Range("E1").Value = "0.0986"
$scratch.Range("F1").Value = "005417"
$scratch.Range("G1").Value = "0.56"
$scratch.Range("A2:C2").NumberFormat = "@"
$scratch.Range("A2").Value = "37.81"
$scratch.Range("B2").Value = "1.24"
$scratch.Range("C2").Value = "0.0069"

# ---------------------------------------------------------------------
# 1. Add the "2022-Q1" worksheet right after "2021-Q4" (i.e. right
#    before "总计", preserving chronological tab order).
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Copy header (B1:H1) + index-column (A2:A3) formatting from the
# "2021-Q4" sheet, which already has the identical column layout and
# the bold/bordered "header / index" style we need to reproduce.
$afterSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$afterSheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$excel.Application.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$scratch.Range("A1").Copy()
$newSheet.Range("B2").PasteSpecial(-4163)
$newSheet.Range("C2").Value = "鹏华尊惠18个月定期开放混合A"
$scratch.Range("B1").Copy()
$newSheet.Range("D2").PasteSpecial(-4163)
$scratch.Range("C1").Copy()
$newSheet.Range("E2").PasteSpecial(-4163)
$scratch.Range("D1").Copy()
$newSheet.Range("F2").PasteSpecial(-4163)
$scratch.Range("E1").Copy()
$newSheet.Range("G2").PasteSpecial(-4163)
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$scratch.Range("F1").Copy()
$newSheet.Range("B3").PasteSpecial(-4163)
$newSheet.Range("C3").Value = "鹏华尊惠18个月定期开放混合C"
$scratch.Range("G1").Copy()
$newSheet.Range("D3").PasteSpecial(-4163)
$scratch.Range("A2").Copy()
$newSheet.Range("E3").PasteSpecial(-4163)
$scratch.Range("B2").Copy()
$newSheet.Range("F3").PasteSpecial(-4163)
$scratch.Range("C2").Copy()
$newSheet.Range("G3").PasteSpecial(-4163)
$newSheet.Range("H3").Value = 8

$excel.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Insert a new top data-row into "总计" for the 2022-Q1 totals,
#    pushing the older quarters (2021-Q4 … 2020-Q4) down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-apply the exact per-column formatting (copied from the row that
# just got pushed down to row 3, which still carries the original
# look: bold/bordered index cell in column A, plain cells in B:D).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)
$excel.Application.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.11

# Renumber the running index in column A for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# ---------------------------------------------------------------------
# Clean up the scratch sheet - it must not remain in the workbook.
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$scratch.Delete()
$excel.DisplayAlerts = $true
